$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value2 = 0.302
$ws.Cells.Item(4, 5).Value2 = 0.16
$ws.Cells.Item(4, 7).Value2 = 0.112
$ws.Cells.Item(4, 8).Value2 = 0.2
$ws.Cells.Item(4, 10).Value2 = 0.101
$ws.Cells.Item(4, 11).Value2 = 0.352
$ws.Cells.Item(4, 12).Value2 = 0.105
$ws.Cells.Item(4, 13).Value2 = 0.323
$ws.Cells.Item(4, 14).Value2 = 0.277
$ws.Cells.Item(4, 15).Value2 = 0.019
$ws.Cells.Item(4, 16).Value2 = 0.139
$ws.Cells.Item(4, 17).Value2 = 0.542
$ws.Cells.Item(4, 18).Value2 = 0.214
$ws.Cells.Item(4, 19).Value2 = 0.462
$ws.Cells.Item(4, 20).Value2 = 0.299
$ws.Cells.Item(4, 23).Value2 = 0.239
$ws.Cells.Item(4, 25).Value2 = 0.212
$ws.Cells.Item(4, 26).Value2 = 0.467
$ws.Cells.Item(4, 27).Value2 = 0.134
$ws.Cells.Item(4, 28).Value2 = 0.365
$ws.Cells.Item(4, 31).Value2 = 0.077
$ws.Cells.Item(4, 32).Value2 = 0.722
$ws.Cells.Item(4, 33).Value2 = 0.095
$ws.Cells.Item(4, 35).Value2 = 0.639
$ws.Cells.Item(4, 36).Value2 = 0.175
$ws.Cells.Item(4, 37).Value2 = 0.419
$ws.Cells.Item(4, 38).Value2 = 0.6860000000000001
$ws.Cells.Item(4, 41).Value2 = 0.6820000000000001

# Row 5
$ws.Cells.Item(5, 2).Value2 = 0.806
$ws.Cells.Item(5, 3).Value2 = 0.157
$ws.Cells.Item(5, 4).Value2 = 0.396
$ws.Cells.Item(5, 5).Value2 = 0.694
$ws.Cells.Item(5, 6).Value2 = 0.212
$ws.Cells.Item(5, 7).Value2 = 0.461
$ws.Cells.Item(5, 8).Value2 = 0.861
$ws.Cells.Item(5, 9).Value2 = 0.12
$ws.Cells.Item(5, 10).Value2 = 0.346
$ws.Cells.Item(5, 11).Value2 = 0.667
$ws.Cells.Item(5, 12).Value2 = 0.222
$ws.Cells.Item(5, 13).Value2 = 0.471
$ws.Cells.Item(5, 14).Value2 = 0.861
$ws.Cells.Item(5, 15).Value2 = 0.12
$ws.Cells.Item(5, 16).Value2 = 0.346
$ws.Cells.Item(5, 17).Value2 = 0.611
$ws.Cells.Item(5, 18).Value2 = 0.238
$ws.Cells.Item(5, 19).Value2 = 0.487
$ws.Cells.Item(5, 20).Value2 = 0.611
$ws.Cells.Item(5, 21).Value2 = 0.238
$ws.Cells.Item(5, 22).Value2 = 0.487
$ws.Cells.Item(5, 23).Value2 = 0.722
$ws.Cells.Item(5, 24).Value2 = 0.201
$ws.Cells.Item(5, 25).Value2 = 0.448
$ws.Cells.Item(5, 26).Value2 = 0.833
$ws.Cells.Item(5, 27).Value2 = 0.139
$ws.Cells.Item(5, 28).Value2 = 0.373
$ws.Cells.Item(5, 29).Value2 = 0.778
$ws.Cells.Item(5, 30).Value2 = 0.173
$ws.Cells.Item(5, 31).Value2 = 0.416
$ws.Cells.Item(5, 32).Value2 = 0.972
$ws.Cells.Item(5, 33).Value2 = 0.027
$ws.Cells.Item(5, 34).Value2 = 0.164
$ws.Cells.Item(5, 35).Value2 = 0.75
$ws.Cells.Item(5, 36).Value2 = 0.188
$ws.Cells.Item(5, 37).Value2 = 0.433
$ws.Cells.Item(5, 38).Value2 = 0.917
$ws.Cells.Item(5, 39).Value2 = 0.076
$ws.Cells.Item(5, 40).Value2 = 0.276
$ws.Cells.Item(5, 41).Value2 = 0.88

# Row 6
$ws.Cells.Item(6, 2).Value2 = 0.439
$ws.Cells.Item(6, 5).Value2 = 0.26
$ws.Cells.Item(6, 8).Value2 = 0.325
$ws.Cells.Item(6, 11).Value2 = 0.461
$ws.Cells.Item(6, 14).Value2 = 0.419
$ws.Cells.Item(6, 17).Value2 = 0.574
$ws.Cells.Item(6, 20).Value2 = 0.402
$ws.Cells.Item(6, 23).Value2 = 0.359
$ws.Cells.Item(6, 26).Value2 = 0.598
$ws.Cells.Item(6, 32).Value2 = 0.829
$ws.Cells.Item(6, 35).Value2 = 0.6899999999999999
$ws.Cells.Item(6, 38).Value2 = 0.785
$ws.Cells.Item(6, 41).Value2 = 0.768

# Row 7
$ws.Cells.Item(7, 2).Value2 = 0.604
$ws.Cells.Item(7, 5).Value2 = 0.416
$ws.Cells.Item(7, 8).Value2 = 0.518
$ws.Cells.Item(7, 11).Value2 = 0.5659999999999999
$ws.Cells.Item(7, 14).Value2 = 0.606
$ws.Cells.Item(7, 17).Value2 = 0.596
$ws.Cells.Item(7, 20).Value2 = 0.506
$ws.Cells.Item(7, 23).Value2 = 0.514
$ws.Cells.Item(7, 26).Value2 = 0.72
$ws.Cells.Item(7, 29).Value2 = 0.388
$ws.Cells.Item(7, 32).Value2 = 0.909
$ws.Cells.Item(7, 35).Value2 = 0.725
$ws.Cells.Item(7, 38).Value2 = 0.859
$ws.Cells.Item(7, 41).Value2 = 0.831

# Row 8
$ws.Cells.Item(8, 2).Value2 = 0.761
$ws.Cells.Item(8, 3).Value2 = 0.156
$ws.Cells.Item(8, 4).Value2 = 0.394
$ws.Cells.Item(8, 5).Value2 = 0.578
$ws.Cells.Item(8, 8).Value2 = 0.742
$ws.Cells.Item(8, 9).Value2 = 0.131
$ws.Cells.Item(8, 10).Value2 = 0.363
$ws.Cells.Item(8, 11).Value2 = 0.591
$ws.Cells.Item(8, 12).Value2 = 0.201
$ws.Cells.Item(8, 13).Value2 = 0.448
$ws.Cells.Item(8, 14).Value2 = 0.777
$ws.Cells.Item(8, 15).Value2 = 0.127
$ws.Cells.Item(8, 16).Value2 = 0.356
$ws.Cells.Item(8, 17).Value2 = 0.58
$ws.Cells.Item(8, 18).Value2 = 0.224
$ws.Cells.Item(8, 19).Value2 = 0.473
$ws.Cells.Item(8, 20).Value2 = 0.528
$ws.Cells.Item(8, 21).Value2 = 0.202
$ws.Cells.Item(8, 22).Value2 = 0.45
$ws.Cells.Item(8, 23).Value2 = 0.653
$ws.Cells.Item(8, 24).Value2 = 0.187
$ws.Cells.Item(8, 25).Value2 = 0.432
$ws.Cells.Item(8, 26).Value2 = 0.765
$ws.Cells.Item(8, 27).Value2 = 0.14
$ws.Cells.Item(8, 28).Value2 = 0.375
$ws.Cells.Item(8, 29).Value2 = 0.674
$ws.Cells.Item(8, 30).Value2 = 0.174
$ws.Cells.Item(8, 31).Value2 = 0.417
$ws.Cells.Item(8, 32).Value2 = 0.887
$ws.Cells.Item(8, 33).Value2 = 0.048
$ws.Cells.Item(8, 34).Value2 = 0.22
$ws.Cells.Item(8, 35).Value2 = 0.74
$ws.Cells.Item(8, 36).Value2 = 0.186
$ws.Cells.Item(8, 37).Value2 = 0.431
$ws.Cells.Item(8, 38).Value2 = 0.886
$ws.Cells.Item(8, 39).Value2 = 0.082
$ws.Cells.Item(8, 40).Value2 = 0.286
$ws.Cells.Item(8, 41).Value2 = 0.838

# Row 9
$ws.Cells.Item(9, 2).Value2 = 0.694
$ws.Cells.Item(9, 3).Value2 = 0.212
$ws.Cells.Item(9, 4).Value2 = 0.461
$ws.Cells.Item(9, 5).Value2 = 0.444
$ws.Cells.Item(9, 6).Value2 = 0.247
$ws.Cells.Item(9, 7).Value2 = 0.497
$ws.Cells.Item(9, 8).Value2 = 0.611
$ws.Cells.Item(9, 9).Value2 = 0.238
$ws.Cells.Item(9, 10).Value2 = 0.487
$ws.Cells.Item(9, 11).Value2 = 0.5
$ws.Cells.Item(9, 14).Value2 = 0.667
$ws.Cells.Item(9, 15).Value2 = 0.222
$ws.Cells.Item(9, 16).Value2 = 0.471
$ws.Cells.Item(9, 17).Value2 = 0.528
$ws.Cells.Item(9, 18).Value2 = 0.249
$ws.Cells.Item(9, 19).Value2 = 0.499
$ws.Cells.Item(9, 20).Value2 = 0.417
$ws.Cells.Item(9, 21).Value2 = 0.243
$ws.Cells.Item(9, 22).Value2 = 0.493
$ws.Cells.Item(9, 23).Value2 = 0.556
$ws.Cells.Item(9, 24).Value2 = 0.247
$ws.Cells.Item(9, 25).Value2 = 0.497
$ws.Cells.Item(9, 26).Value2 = 0.667
$ws.Cells.Item(9, 27).Value2 = 0.222
$ws.Cells.Item(9, 28).Value2 = 0.471
$ws.Cells.Item(9, 29).Value2 = 0.583
$ws.Cells.Item(9, 30).Value2 = 0.243
$ws.Cells.Item(9, 31).Value2 = 0.493
$ws.Cells.Item(9, 32).Value2 = 0.75
$ws.Cells.Item(9, 33).Value2 = 0.188
$ws.Cells.Item(9, 34).Value2 = 0.433
$ws.Cells.Item(9, 35).Value2 = 0.722
$ws.Cells.Item(9, 36).Value2 = 0.201
$ws.Cells.Item(9, 37).Value2 = 0.448
$ws.Cells.Item(9, 38).Value2 = 0.833
$ws.Cells.Item(9, 39).Value2 = 0.139
$ws.Cells.Item(9, 40).Value2 = 0.373
$ws.Cells.Item(9, 41).Value2 = 0.768

# Row 10
$ws.Cells.Item(10, 2).Value2 = 0.806
$ws.Cells.Item(10, 3).Value2 = 0.157
$ws.Cells.Item(10, 4).Value2 = 0.396
$ws.Cells.Item(10, 5).Value2 = 0.611
$ws.Cells.Item(10, 6).Value2 = 0.238
$ws.Cells.Item(10, 7).Value2 = 0.487
$ws.Cells.Item(10, 8).Value2 = 0.778
$ws.Cells.Item(10, 9).Value2 = 0.173
$ws.Cells.Item(10, 10).Value2 = 0.416
$ws.Cells.Item(10, 11).Value2 = 0.667
$ws.Cells.Item(10, 12).Value2 = 0.222
$ws.Cells.Item(10, 13).Value2 = 0.471
$ws.Cells.Item(10, 14).Value2 = 0.833
$ws.Cells.Item(10, 15).Value2 = 0.139
$ws.Cells.Item(10, 16).Value2 = 0.373
$ws.Cells.Item(10, 17).Value2 = 0.611
$ws.Cells.Item(10, 18).Value2 = 0.238
$ws.Cells.Item(10, 19).Value2 = 0.487
$ws.Cells.Item(10, 20).Value2 = 0.611
$ws.Cells.Item(10, 21).Value2 = 0.238
$ws.Cells.Item(10, 22).Value2 = 0.487
$ws.Cells.Item(10, 23).Value2 = 0.722
$ws.Cells.Item(10, 24).Value2 = 0.201
$ws.Cells.Item(10, 25).Value2 = 0.448
$ws.Cells.Item(10, 26).Value2 = 0.833
$ws.Cells.Item(10, 27).Value2 = 0.139
$ws.Cells.Item(10, 28).Value2 = 0.373
$ws.Cells.Item(10, 29).Value2 = 0.667
$ws.Cells.Item(10, 30).Value2 = 0.222
$ws.Cells.Item(10, 31).Value2 = 0.471
$ws.Cells.Item(10, 32).Value2 = 0.972
$ws.Cells.Item(10, 33).Value2 = 0.027
$ws.Cells.Item(10, 34).Value2 = 0.164
$ws.Cells.Item(10, 35).Value2 = 0.75
$ws.Cells.Item(10, 36).Value2 = 0.188
$ws.Cells.Item(10, 37).Value2 = 0.433
$ws.Cells.Item(10, 38).Value2 = 0.917
$ws.Cells.Item(10, 39).Value2 = 0.076
$ws.Cells.Item(10, 40).Value2 = 0.276
$ws.Cells.Item(10, 41).Value2 = 0.88

# Row 11
$ws.Cells.Item(11, 2).Value2 = 0.806
$ws.Cells.Item(11, 3).Value2 = 0.157
$ws.Cells.Item(11, 4).Value2 = 0.396
$ws.Cells.Item(11, 5).Value2 = 0.694
$ws.Cells.Item(11, 6).Value2 = 0.212
$ws.Cells.Item(11, 7).Value2 = 0.461
$ws.Cells.Item(11, 8).Value2 = 0.861
$ws.Cells.Item(11, 9).Value2 = 0.12
$ws.Cells.Item(11, 10).Value2 = 0.346
$ws.Cells.Item(11, 11).Value2 = 0.667
$ws.Cells.Item(11, 12).Value2 = 0.222
$ws.Cells.Item(11, 13).Value2 = 0.471
$ws.Cells.Item(11, 14).Value2 = 0.861
$ws.Cells.Item(11, 15).Value2 = 0.12
$ws.Cells.Item(11, 16).Value2 = 0.346
$ws.Cells.Item(11, 17).Value2 = 0.611
$ws.Cells.Item(11, 18).Value2 = 0.238
$ws.Cells.Item(11, 19).Value2 = 0.487
$ws.Cells.Item(11, 20).Value2 = 0.611
$ws.Cells.Item(11, 21).Value2 = 0.238
$ws.Cells.Item(11, 22).Value2 = 0.487
$ws.Cells.Item(11, 23).Value2 = 0.722
$ws.Cells.Item(11, 24).Value2 = 0.201
$ws.Cells.Item(11, 25).Value2 = 0.448
$ws.Cells.Item(11, 26).Value2 = 0.833
$ws.Cells.Item(11, 27).Value2 = 0.139
$ws.Cells.Item(11, 28).Value2 = 0.373
$ws.Cells.Item(11, 29).Value2 = 0.722
$ws.Cells.Item(11, 30).Value2 = 0.201
$ws.Cells.Item(11, 31).Value2 = 0.448
$ws.Cells.Item(11, 32).Value2 = 0.972
$ws.Cells.Item(11, 33).Value2 = 0.027
$ws.Cells.Item(11, 34).Value2 = 0.164
$ws.Cells.Item(11, 35).Value2 = 0.75
$ws.Cells.Item(11, 36).Value2 = 0.188
$ws.Cells.Item(11, 37).Value2 = 0.433
$ws.Cells.Item(11, 38).Value2 = 0.917
$ws.Cells.Item(11, 39).Value2 = 0.076
$ws.Cells.Item(11, 40).Value2 = 0.276
$ws.Cells.Item(11, 41).Value2 = 0.88

# Row 12
$ws.Cells.Item(12, 2).Value2 = 1.172
$ws.Cells.Item(12, 3).Value2 = 0.212
$ws.Cells.Item(12, 4).Value2 = 0.46
$ws.Cells.Item(12, 5).Value2 = 1.68
$ws.Cells.Item(12, 6).Value2 = 1.098
$ws.Cells.Item(12, 7).Value2 = 1.048
$ws.Cells.Item(12, 8).Value2 = 1.613
$ws.Cells.Item(12, 9).Value2 = 1.334
$ws.Cells.Item(12, 10).Value2 = 1.155
$ws.Cells.Item(12, 11).Value2 = 1.417
$ws.Cells.Item(12, 12).Value2 = 0.576
$ws.Cells.Item(12, 13).Value2 = 0.759
$ws.Cells.Item(12, 14).Value2 = 1.355
$ws.Cells.Item(12, 15).Value2 = 0.552
$ws.Cells.Item(12, 16).Value2 = 0.743
$ws.Cells.Item(12, 26).Value2 = 1.267
$ws.Cells.Item(12, 27).Value2 = 0.329
$ws.Cells.Item(12, 28).Value2 = 0.573
$ws.Cells.Item(12, 29).Value2 = 1.786
$ws.Cells.Item(12, 30).Value2 = 2.526
$ws.Cells.Item(12, 31).Value2 = 1.589
$ws.Cells.Item(12, 32).Value2 = 1.257
$ws.Cells.Item(12, 33).Value2 = 0.248
$ws.Cells.Item(12, 34).Value2 = 0.498
$ws.Cells.Item(12, 35).Value2 = 1.037
$ws.Cells.Item(12, 36).Value2 = 0.036
$ws.Cells.Item(12, 37).Value2 = 0.189
$ws.Cells.Item(12, 38).Value2 = 1.091
$ws.Cells.Item(12, 39).Value2 = 0.083
$ws.Cells.Item(12, 40).Value2 = 0.287
$ws.Cells.Item(12, 41).Value2 = 1.128

# Row 13
$ws.Cells.Item(13, 2).Value2 = 3.389
$ws.Cells.Item(13, 3).Value2 = 1.404
$ws.Cells.Item(13, 4).Value2 = 1.185
$ws.Cells.Item(13, 5).Value2 = 4.594
$ws.Cells.Item(13, 6).Value2 = 0.429
$ws.Cells.Item(13, 7).Value2 = 0.655
$ws.Cells.Item(13, 8).Value2 = 4.611
$ws.Cells.Item(13, 9).Value2 = 0.627
$ws.Cells.Item(13, 10).Value2 = 0.792
$ws.Cells.Item(13, 11).Value2 = 2.281
$ws.Cells.Item(13, 12).Value2 = 0.577
$ws.Cells.Item(13, 13).Value2 = 0.76
$ws.Cells.Item(13, 14).Value2 = 3.25
$ws.Cells.Item(13, 15).Value2 = 0.743
$ws.Cells.Item(13, 16).Value2 = 0.862
$ws.Cells.Item(13, 26).Value2 = 2.5
$ws.Cells.Item(13, 27).Value2 = 2.956
$ws.Cells.Item(13, 28).Value2 = 1.719
$ws.Cells.Item(13, 29).Value2 = 6.314
$ws.Cells.Item(13, 30).Value2 = 2.216
$ws.Cells.Item(13, 31).Value2 = 1.488
$ws.Cells.Item(13, 32).Value2 = 1.639
$ws.Cells.Item(13, 33).Value2 = 0.731
$ws.Cells.Item(13, 34).Value2 = 0.855
$ws.Cells.Item(13, 35).Value2 = 1.306
$ws.Cells.Item(13, 36).Value2 = 0.379
$ws.Cells.Item(13, 37).Value2 = 0.616
$ws.Cells.Item(13, 38).Value2 = 1.611
$ws.Cells.Item(13, 39).Value2 = 0.738
$ws.Cells.Item(13, 40).Value2 = 0.859
$ws.Cells.Item(13, 41).Value2 = 1.519
